# Applies the "Updated symbol list" commit: refreshes prices/volume percentages
# for existing rows, and for a handful of rows the underlying coin (name + link)
# was swapped out along with its price/volume while staying on the same row.
#
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora
# D/E hold numeric-looking text (e.g. "247.01", "0.86%") that must stay stored
# as literal text (matching the original inlineStr cells), not be coerced into
# numbers/percentages by Excel's automatic type detection. We force text by
# temporarily applying a "@" (Text) number format before assigning the value,
# then reverting the cell style to "Normal" so no stray style index lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text fields (coin name / link) -- safe to assign directly.
$textCells = @{
    "B9" = "WazirX"
    "C9" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B20" = "LiechtensteinCryptoassetsExchange"
    "C20" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}

# Numeric-looking fields (price / percentage) -- force text storage.
$numericTextCells = @{
    "D2" = "247.01"
    "E2" = "0.86%"
    "E3" = "5.14%"
    "D4" = "5.083"
    "E4" = "2.74%"
    "D5" = "0.05602"
    "E5" = "-0.25%"
    "D6" = "6.491"
    "E6" = "-0.82%"
    "D7" = "0.8134"
    "E7" = "0.26%"
    "D8" = "0.8447"
    "E8" = "0.52%"
    "D9" = "0.1338"
    "E9" = "0.14%"
    "D10" = "0.06980"
    "E10" = "0.20%"
    "D11" = "0.02851"
    "E11" = "0.09%"
    "D12" = "0.09388"
    "E12" = "-0.24%"
    "D13" = "0.001510"
    "E13" = "-1.05%"
    "D14" = "0.0005959"
    "E14" = "-93.87%"
    "D15" = "0.006150"
    "E15" = "-1.21%"
    "D16" = "3.612"
    "E16" = "3.19%"
    "E17" = "0.32%"
    "E18" = "-1.71%"
    "E19" = "-1.31%"
    "D20" = "0.03186"
    "E20" = "-1.70%"
    "D21" = "0.1298"
    "E21" = "0.50%"
    "D22" = "3.743"
    "E22" = "0.15%"
    "D23" = "0.04656"
    "E23" = "-0.69%"
    "E24" = "-1.44%"
    "D25" = "0.001248"
    "E25" = "0.46%"
    "E26" = "1.38%"
    "D27" = "0.00009596"
    "E27" = "-1.06%"
    "E28" = "-27.95%"
    "E40" = "0.55%"
    "D41" = "0.006157"
    "E41" = "-1.15%"
    "D42" = "0.1059"
    "E42" = "-21.69%"
    "D43" = "0.002499"
    "E43" = "-8.23%"
    "D44" = "0.008935"
    "E44" = "10.67%"
    "D45" = "0.00005299"
    "E45" = "0.46%"
    "E47" = "-38.88%"
    "D48" = "0.002617"
    "E48" = "28.08%"
}

foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
    $cell.Style = "Normal"
}
